# Applies the BankDeductDtl.xlsx edit described in the commit diff:
#  - Rewords several "CdCode:" reference notes (in column G of sheet "DBD")
#    from a colon-and-period style to a dot-and-colon style, and reflows
#    a couple of multi-line notes onto more lines.
#  - Adjusts the row heights of rows 24 and 38 on "DBD" to fit the
#    rewrapped text.
#  - Updates the saved selection / scroll position on "DBD" to reflect
#    where the editor was last working (cell G40, scrolled so row 34 is
#    the top visible row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

# --- Column G reference-note rewording -------------------------------

$ws.Range("G12").Value = "CdCode.RepayType`n1:期款`n2:部分償還`n3:結案`n4:帳管費`n5:火險費`n6:契變手續費`n7:法務費`n9:其他"

$ws.Range("G16").Value = "L4451建檔交易產生者，此欄位由額度檔抓取`nCdCode.BankCdAppl"

$ws.Range("G24").Value = "L4451建檔交易產生者，此欄位由額度檔抓取`nCdCode.PostDepCode`nG:劃撥`nP:存簿"

$ws.Range("G26").Value = "L4451建檔交易產生者，此欄位由額度檔抓取`nCdCode.RelationCode"

$ws.Range("G30").Value = "CdCode.Sex"

$ws.Range("G37").Value = "CdCode.AmlCheckItem`n0:非可疑名單/已完成名單確認`n1:需審查/確認`n2:為凍結名單/未確定名單"

$ws.Range("G38").Value = "空白:未回`n00:扣款成功`n>00:扣款失敗`n失敗原因 ref. CdCode.ProcCode 處理說明`n ACH ：002 + ReturnCode(2)`n 郵局：003 + ReturnCode(2)"

# --- Row heights, resized to fit the rewrapped notes above ------------

$ws.Rows.Item(24).RowHeight = 64.8
$ws.Rows.Item(38).RowHeight = 113.4

# --- Saved view state: active selection --------------------------------
# (the workbook was also scrolled so row 34 was the top visible row, but
# that window-scroll state isn't reachable through this object model)

$ws.Range("G40").Select()
